$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Sheet, $Address, $Text) {
    $cell = $Sheet.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.ClearFormats()
}

Set-TextValue $ws 'D2' '246.07'
Set-TextValue $ws 'D3' '22.78'
Set-TextValue $ws 'D4' '5.426'
Set-TextValue $ws 'D5' '0.05743'
Set-TextValue $ws 'D6' '3.431'
Set-TextValue $ws 'B7' 'KuCoinToken'
Set-TextValue $ws 'C7' 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
Set-TextValue $ws 'D7' '6.330'
Set-TextValue $ws 'E7' '6KuCoinTokenKCS'
Set-TextValue $ws 'B8' 'MXToken'
Set-TextValue $ws 'C8' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws 'D8' '0.8138'
Set-TextValue $ws 'E8' '7MXTokenMX'
Set-TextValue $ws 'B9' 'FTXToken'
Set-TextValue $ws 'C9' 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue $ws 'D9' '0.8845'
Set-TextValue $ws 'E9' '8FTXTokenFTT'
Set-TextValue $ws 'B10' 'WazirX'
Set-TextValue $ws 'C10' 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue $ws 'D10' '0.1442'
Set-TextValue $ws 'E10' '9WazirXWRX'
Set-TextValue $ws 'B11' 'MandalaExchangeToken'
Set-TextValue $ws 'C11' 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue $ws 'D11' '0.07343'
Set-TextValue $ws 'E11' '10MandalaExchangeTokenMDX'
Set-TextValue $ws 'B12' 'LiechtensteinCryptoassetsExchange'
Set-TextValue $ws 'C12' 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue $ws 'D12' '0.02990'
Set-TextValue $ws 'E12' '11LiechtensteinCryptoassetsExchangeLCX'
Set-TextValue $ws 'B13' 'BitrueCoin'
Set-TextValue $ws 'C13' 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue $ws 'D13' '0.03061'
Set-TextValue $ws 'E13' '12BitrueCoinBTR'
Set-TextValue $ws 'B14' 'BitMartToken'
Set-TextValue $ws 'C14' 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue $ws 'D14' '0.09391'
Set-TextValue $ws 'E14' '13BitMartTokenBMX'
Set-TextValue $ws 'B15' 'BitForexToken'
Set-TextValue $ws 'C15' 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue $ws 'D15' '0.001579'
Set-TextValue $ws 'E15' '14BitForexTokenBF'
Set-TextValue $ws 'B16' 'CoinExToken'
Set-TextValue $ws 'C16' 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
Set-TextValue $ws 'D16' '0.04821'
Set-TextValue $ws 'E16' '15CoinExTokenCET'
Set-TextValue $ws 'B17' 'One'
Set-TextValue $ws 'C17' 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextValue $ws 'D17' '0.0005840'
Set-TextValue $ws 'E17' '16OneONE'
Set-TextValue $ws 'B18' 'TigerCash'
Set-TextValue $ws 'C18' 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue $ws 'D18' '0.006160'
Set-TextValue $ws 'E18' '17TigerCashTCH'
Set-TextValue $ws 'B19' 'HotbitToken'
Set-TextValue $ws 'C19' 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
Set-TextValue $ws 'D19' '0.005107'
Set-TextValue $ws 'E19' '18HotbitTokenHTB'
Set-TextValue $ws 'B20' 'BitKan'
Set-TextValue $ws 'C20' 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
Set-TextValue $ws 'D20' '0.0009969'
Set-TextValue $ws 'E20' '19BitKanKAN'
Set-TextValue $ws 'B21' 'NitroEx'
Set-TextValue $ws 'C21' 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
Set-TextValue $ws 'D21' '0.0001500'
Set-TextValue $ws 'E21' '20NitroExNTX'
Set-TextValue $ws 'B22' 'LEO'
Set-TextValue $ws 'C22' 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws 'D22' '3.749'
Set-TextValue $ws 'E22' '21LEOLEO'
Set-TextValue $ws 'B23' 'BTSEToken'
Set-TextValue $ws 'C23' 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextValue $ws 'D23' '2.199'
Set-TextValue $ws 'E23' '22BTSETokenBTSE'
Set-TextValue $ws 'D24' '0.3276'
Set-TextValue $ws 'D25' '0.1316'
Set-TextValue $ws 'D26' '4.169'
Set-TextValue $ws 'D40' '0.03916'
Set-TextValue $ws 'D41' '0.006774'
Set-TextValue $ws 'E41' '40KickTokenKICKBestin24h'
Set-TextValue $ws 'D44' '0.007140'
Set-TextValue $ws 'D45' '0.00005644'
Set-TextValue $ws 'D47' '0.3800'
Set-TextValue $ws 'D48' '0.1685'
